# Auto-generated edit script applying the Asura_Profits.xlsx market-data refresh
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1510.8
$ws.Range("I86").Value = 867.6667
$ws.Range("J86").Value = 1786.4286
$ws.Range("K86").Value = 867.6667
$ws.Range("L86").Value = 1786.4286
$ws.Range("M86").Value = 255.3333
$ws.Range("N86").Value = -4032.4286
$ws.Range("H89").Value = 1510.8
$ws.Range("I89").Value = 867.6667
$ws.Range("J89").Value = 1786.4286
$ws.Range("K89").Value = 4338.3335
$ws.Range("L89").Value = 8932.143
$ws.Range("M89").Value = 1277.6665
$ws.Range("N89").Value = -20164.143
$ws.Range("H106").Value = 4335
$ws.Range("I106").Value = 5502.5
$ws.Range("K106").Value = 5502.5
$ws.Range("M106").Value = -4871.5
$ws.Range("H112").Value = 5146.615
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 5146.615
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 15439.845
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -17655.845
$ws.Range("H125").Value = 3918.6667
$ws.Range("I125").Value = 1766
$ws.Range("J125").Value = 4533.7144
$ws.Range("K125").Value = 15894
$ws.Range("L125").Value = 40803.4296
$ws.Range("M125").Value = -13434
$ws.Range("N125").Value = -45723.4296
$ws.Range("H129").Value = 1118.1305
$ws.Range("J129").Value = 1128.8383
$ws.Range("L129").Value = 3386.5149
$ws.Range("N129").Value = -13386.5149
$ws.Range("H132").Value = 2558.074
$ws.Range("J132").Value = 2109.077
$ws.Range("L132").Value = 6327.231000000001
$ws.Range("N132").Value = -11387.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2920.9524
$ws.Range("I61").Value = 3357
$ws.Range("K61").Value = 3357
$ws.Range("M61").Value = -3145
$ws.Range("H97").Value = 738.0454999999999
$ws.Range("I97").Value = 738.0454999999999
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 738.0454999999999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -242.0454999999999
$ws.Range("N97").ClearContents()
$ws.Range("H136").Value = 2920.9524
$ws.Range("I136").Value = 3357
$ws.Range("K136").Value = 10071
$ws.Range("M136").Value = -7521

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 36956
$ws.Range("J81").Value = 36956
$ws.Range("L81").Value = 36956
$ws.Range("N81").Value = -39078
$ws.Range("H84").Value = 36956
$ws.Range("J84").Value = 36956
$ws.Range("L84").Value = 110868
$ws.Range("N84").Value = -121476
$ws.Range("H105").Value = 5716851.5
$ws.Range("I105").Value = 7938939.5
$ws.Range("J105").Value = 2911
$ws.Range("K105").Value = 7938939.5
$ws.Range("L105").Value = 2911
$ws.Range("M105").Value = -7937192.5
$ws.Range("N105").Value = -6405
$ws.Range("H107").Value = 25978.818
$ws.Range("I107").Value = 36913.867
$ws.Range("J107").Value = 2546.5715
$ws.Range("K107").Value = 36913.867
$ws.Range("L107").Value = 2546.5715
$ws.Range("M107").Value = -34993.867
$ws.Range("N107").Value = -6386.5715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 31649.875
$ws.Range("J52").Value = 31649.875
$ws.Range("L52").Value = 31649.875
$ws.Range("N52").Value = -32237.875
$ws.Range("H62").Value = 44433.75
$ws.Range("I62").Value = 52720.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 52720.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -52096.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 44433.75
$ws.Range("I65").Value = 52720.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 263602.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -260482.5
$ws.Range("N65").Value = -21240
$ws.Range("H132").Value = 425711.03
$ws.Range("I132").Value = 617243.4399999999
$ws.Range("J132").Value = 4339.8
$ws.Range("K132").Value = 1851730.32
$ws.Range("L132").Value = 13019.4
$ws.Range("M132").Value = -1849200.32
$ws.Range("N132").Value = -18079.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 150124.56
$ws.Range("I68").Value = 270864.9
$ws.Range("K68").Value = 812594.7000000001
$ws.Range("M68").Value = -811783.7000000001
$ws.Range("H71").Value = 150124.56
$ws.Range("I71").Value = 270864.9
$ws.Range("K71").Value = 2437784.1
$ws.Range("M71").Value = -2433728.1
$ws.Range("H107").Value = 1122.127
$ws.Range("I107").Value = 726.8542
$ws.Range("J107").Value = 2387
$ws.Range("K107").Value = 2180.5626
$ws.Range("L107").Value = 7161
$ws.Range("M107").Value = -260.5626000000002
$ws.Range("N107").Value = -11001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 41638
$ws.Range("J51").Value = 41638
$ws.Range("L51").Value = 41638
$ws.Range("N51").Value = -42656
$ws.Range("H70").Value = 7604.16
$ws.Range("I70").Value = 7700.174
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 7700.174
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -7430.174
$ws.Range("N70").Value = -7040
$ws.Range("H73").Value = 7604.16
$ws.Range("I73").Value = 7700.174
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 7700.174
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -6764.174
$ws.Range("N73").Value = -8372
$ws.Range("H126").Value = 4188.778
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 4814.2856
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 14442.8568
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -19382.8568
$ws.Range("H132").Value = 3204.087
$ws.Range("I132").Value = 2363
$ws.Range("J132").Value = 4512.4443
$ws.Range("K132").Value = 7089
$ws.Range("L132").Value = 13537.3329
$ws.Range("M132").Value = -4559
$ws.Range("N132").Value = -18597.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5523.2915
$ws.Range("I132").Value = 5545.3687
$ws.Range("J132").Value = 5439.4
$ws.Range("K132").Value = 16636.1061
$ws.Range("L132").Value = 16318.2
$ws.Range("M132").Value = -14106.1061
$ws.Range("N132").Value = -21378.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 10001816
$ws.Range("I122").Value = 16668266
$ws.Range("J122").Value = 2141
$ws.Range("K122").Value = 50004798
$ws.Range("L122").Value = 6423
$ws.Range("M122").Value = -50002348
$ws.Range("N122").Value = -11323
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
